$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D (old D:K shifts to F:M)
$ws.Columns("D:E").Insert()

# Copy number formatting (date / number styles) from column F into the new D:E columns
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new columns D (newest quarter) and E with the reported figures
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 154200
$ws.Range("E8").Value = 173100
$ws.Range("D9").Value = 79000
$ws.Range("E9").Value = 87500
$ws.Range("D10").Value = 75200
$ws.Range("E10").Value = 85600
$ws.Range("D12").Value = 20700
$ws.Range("E12").Value = 18500
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 3800
$ws.Range("E14").Value = 400
$ws.Range("D15").Value = 4100
$ws.Range("E15").Value = 3600
$ws.Range("D17").Value = 134600
$ws.Range("E17").Value = 133200
$ws.Range("D18").Value = 19600
$ws.Range("E18").Value = 39900
$ws.Range("D20").Value = 800
$ws.Range("E20").Value = 400
$ws.Range("D21").Value = 24600
$ws.Range("E21").Value = 43800
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 20500
$ws.Range("E23").Value = 40300
$ws.Range("D24").Value = -200
$ws.Range("E24").Value = 2700
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 20700
$ws.Range("E26").Value = 37600
$ws.Range("D27").Value = 20700
$ws.Range("E27").Value = 37500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -1300
$ws.Range("E29").Value = -2800
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -800
$ws.Range("E32").Value = -400
$ws.Range("D33").Value = 19400
$ws.Range("E33").Value = 34800
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 19400
$ws.Range("E35").Value = 34800
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 349300
$ws.Range("E41").Value = 338700
$ws.Range("D42").Value = 2500
$ws.Range("E42").Value = 3100
$ws.Range("D43").Value = 102700
$ws.Range("E43").Value = 114700
$ws.Range("D44").Value = 98000
$ws.Range("E44").Value = 110300
$ws.Range("D45").Value = 16000
$ws.Range("E45").Value = 18100
$ws.Range("D46").Value = 568400
$ws.Range("E46").Value = 584800
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 31300
$ws.Range("E48").Value = 30200
$ws.Range("D49").Value = 156800
$ws.Range("E49").Value = 157900
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 60000
$ws.Range("E52").Value = 60800
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 816500
$ws.Range("E54").Value = 833600
$ws.Range("D57").Value = 39600
$ws.Range("E57").Value = 45600
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 70700
$ws.Range("E59").Value = 79000
$ws.Range("D60").Value = 110300
$ws.Range("E60").Value = 124600
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 98900
$ws.Range("E62").Value = 109500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 209700
$ws.Range("E66").Value = 234600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 512800
$ws.Range("E72").Value = 480600
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 606800
$ws.Range("E76").Value = 599100
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 19400
$ws.Range("E81").Value = 34800
$ws.Range("D83").Value = 4100
$ws.Range("E83").Value = 3600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 39500
$ws.Range("E89").Value = 28100
$ws.Range("D91").Value = -3700
$ws.Range("E91").Value = -7200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -3100
$ws.Range("E94").Value = -91900
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -25500
$ws.Range("E100").Value = -31000
$ws.Range("D101").Value = -300
$ws.Range("E101").Value = 400
$ws.Range("D102").Value = 10600
$ws.Range("E102").Value = -94300
